# The "Income" sheet's entries are reordered/extended:
#   - the old row 2 (Salary, 1000) is replaced and moved to the bottom
#     with an updated amount (8000)
#   - the old row 3 (Interest from saving account, 5000) becomes the new
#     row 2, keeping its amount/date
#   - two new income rows are inserted in between: Lottery (1500) and
#     Gift (1000)
# Column C holds dates (serial numbers) formatted with the existing date
# style already used by the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Interest from saving account "
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 45898.229537037034

$ws.Range("A3").Value = "Lottery"
$ws.Range("B3").Value = 1500
$ws.Range("C3").Value = 45895.229537037034

$ws.Range("A4").Value = "Gift"
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 45891.229537037034

$ws.Range("A5").Value = "Salary"
$ws.Range("B5").Value = 8000
$ws.Range("C5").Value = 45888.229537037034

# Copy the date cell's formatting (built-in m/d/yyyy style already used by
# C2/C3) onto the two newly added date cells, so a duplicate custom number
# format isn't created in the workbook's style table.
$ws.Range("C2").Copy()
$ws.Range("C4:C5").PasteSpecial(-4122)
